$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 40.92730033333333
$ws.Range("N2").Value = 122.781901
$ws.Range("O2").Value = 0.3921621516522625
$ws.Range("P2").Value = 0.3921621516522625
$ws.Range("Q2").Value = 1667.926286701197
$ws.Range("R2").Value = 15011.33658031077
$ws.Range("S2").Value = 0.008339471808535916
$ws.Range("T2").Value = 0.008339471808535915

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.598228333333334
$ws.Range("N3").Value = 22.794685
$ws.Range("O3").Value = 0.07280562235174674
$ws.Range("P3").Value = 0.07280562235174673
$ws.Range("Q3").Value = 309.6535727083545
$ws.Range("R3").Value = 2786.88215437519
$ws.Range("S3").Value = 0.00154823822887346
$ws.Range("T3").Value = 0.00154823822887346

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 49.40125466666666
$ws.Range("N4").Value = 148.203764
$ws.Range("O4").Value = 0.4733589112063359
$ws.Range("P4").Value = 0.4733589112063358
$ws.Range("Q4").Value = 2013.268663788326
$ws.Range("R4").Value = 18119.41797409493
$ws.Range("S4").Value = 0.01006615064379
$ws.Range("T4").Value = 0.01006615064379

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.436425
$ws.Range("N5").Value = 19.309275
$ws.Range("O5").Value = 0.06167331478965488
$ws.Range("P5").Value = 0.06167331478965488
$ws.Range("Q5").Value = 262.30614681265
$ws.Range("R5").Value = 2360.75532131385
$ws.Range("S5").Value = 0.001311505630669192
$ws.Range("T5").Value = 0.001311505630669192

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.92730033333333
$ws.Range("N6").Value = 122.781901
$ws.Range("O6").Value = 0.3921621516522625
$ws.Range("P6").Value = 0.3921621516522625
$ws.Range("Q6").Value = 69138.05080383508
$ws.Range("R6").Value = 622242.4572345158
$ws.Range("S6").Value = 0.3456836373243142
$ws.Range("T6").Value = 0.3456836373243141

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.598228333333334
$ws.Range("N7").Value = 22.794685
$ws.Range("O7").Value = 0.07280562235174674
$ws.Range("P7").Value = 0.07280562235174673
$ws.Range("Q7").Value = 12835.60587311169
$ws.Range("R7").Value = 115520.4528580052
$ws.Range("S7").Value = 0.06417680096402796
$ws.Range("T7").Value = 0.06417680096402795

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 49.40125466666666
$ws.Range("N8").Value = 148.203764
$ws.Range("O8").Value = 0.4733589112063359
$ws.Range("P8").Value = 0.4733589112063358
$ws.Range("Q8").Value = 83453.01124431675
$ws.Range("R8").Value = 751077.1011988507
$ws.Range("S8").Value = 0.4172570695470356
$ws.Range("T8").Value = 0.4172570695470356

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.436425
$ws.Range("N9").Value = 19.309275
$ws.Range("O9").Value = 0.06167331478965488
$ws.Range("P9").Value = 0.06167331478965488
$ws.Range("Q9").Value = 10872.983925662
$ws.Range("R9").Value = 97856.85533095799
$ws.Range("S9").Value = 0.05436387905490603
$ws.Range("T9").Value = 0.05436387905490603

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 40.92730033333333
$ws.Range("N10").Value = 122.781901
$ws.Range("O10").Value = 0.3921621516522625
$ws.Range("P10").Value = 0.3921621516522625
$ws.Range("Q10").Value = 4132.241344572167
$ws.Range("R10").Value = 37190.1721011495
$ws.Range("S10").Value = 0.02066081125640275
$ws.Range("T10").Value = 0.02066081125640275

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 7.598228333333334
$ws.Range("N11").Value = 22.794685
$ws.Range("O11").Value = 0.07280562235174674
$ws.Range("P11").Value = 0.07280562235174673
$ws.Range("Q11").Value = 767.1581806955328
$ws.Range("R11").Value = 6904.423626259795
$ws.Range("S11").Value = 0.003835717484404766
$ws.Range("T11").Value = 0.003835717484404765

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 49.40125466666666
$ws.Range("N12").Value = 148.203764
$ws.Range("O12").Value = 0.4733589112063359
$ws.Range("P12").Value = 0.4733589112063358
$ws.Range("Q12").Value = 4987.81755319146
$ws.Range("R12").Value = 44890.35797872314
$ws.Range("S12").Value = 0.0249386104185865
$ws.Range("T12").Value = 0.0249386104185865

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.436425
$ws.Range("N13").Value = 19.309275
$ws.Range("O13").Value = 0.06167331478965488
$ws.Range("P13").Value = 0.06167331478965488
$ws.Range("Q13").Value = 649.856239713325
$ws.Range("R13").Value = 5848.706157419925
$ws.Range("S13").Value = 0.00324921900560064
$ws.Range("T13").Value = 0.00324921900560064

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 40.92730033333333
$ws.Range("N14").Value = 122.781901
$ws.Range("O14").Value = 0.3921621516522625
$ws.Range("P14").Value = 0.3921621516522625
$ws.Range("Q14").Value = 3495.713162406449
$ws.Range("R14").Value = 31461.41846165804
$ws.Range("S14").Value = 0.01747823126300968
$ws.Range("T14").Value = 0.01747823126300968

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 7.598228333333334
$ws.Range("N15").Value = 22.794685
$ws.Range("O15").Value = 0.07280562235174674
$ws.Range("P15").Value = 0.07280562235174673
$ws.Range("Q15").Value = 648.9855568159745
$ws.Range("R15").Value = 5840.870011343771
$ws.Range("S15").Value = 0.003244865674440551
$ws.Range("T15").Value = 0.003244865674440551

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 49.40125466666666
$ws.Range("N16").Value = 148.203764
$ws.Range("O16").Value = 0.4733589112063359
$ws.Range("P16").Value = 0.4733589112063358
$ws.Range("Q16").Value = 4219.496882793654
$ws.Range("R16").Value = 37975.47194514288
$ws.Range("S16").Value = 0.02109708059692372
$ws.Range("T16").Value = 0.02109708059692372

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.436425
$ws.Range("N17").Value = 19.309275
$ws.Range("O17").Value = 0.06167331478965488
$ws.Range("P17").Value = 0.06167331478965488
$ws.Range("Q17").Value = 549.75274225495
$ws.Range("R17").Value = 4947.77468029455
$ws.Range("S17").Value = 0.002748711098479012
$ws.Range("T17").Value = 0.002748711098479012

